$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (data refreshed for the next day: 2024-09-24 -> 2024-09-25)
$ws.Name = "IClientBalance-20240925-085931-"

# Update the "Dt. Referencia" column (G) for every data row (2..274) from 45559 to 45560
For ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45560
}

# A handful of rows also had their underlying balances refreshed upstream.
# D = Vl. Projetado, E = Saldo Previsto, H = Vl. Total (= D + E)
$ws.Cells.Item(102, 5).Value = 617.67999999999995
$ws.Cells.Item(102, 8).Value = 617.67999999999995

$ws.Cells.Item(105, 5).Value = 565
$ws.Cells.Item(105, 8).Value = 565

$ws.Cells.Item(108, 4).Value = 41280.080000000002
$ws.Cells.Item(108, 8).Value = 44251.58

$ws.Cells.Item(110, 5).Value = 665.04
$ws.Cells.Item(110, 8).Value = 665.04

$ws.Cells.Item(118, 5).Value = 916.85
$ws.Cells.Item(118, 8).Value = 916.85
